$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while guaranteeing it is stored as
# text (not auto-coerced to a Number/Date by the input parser), and
# without leaving a residual numeric-format style behind on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '28.565.33'
Set-TextValue 'E2' '  +1.89%  '
Set-TextValue 'D3' '1.826.61'
Set-TextValue 'E3' '  +1.89%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'D5' '317.50'
Set-TextValue 'E5' '  +0.18%  '
Set-TextValue 'E6' '  +0.06%  '
Set-TextValue 'D7' '0.5411'
Set-TextValue 'E7' '  +1.07%  '
Set-TextValue 'D8' '0.4008'
Set-TextValue 'E8' '  +6.32%  '
Set-TextValue 'D9' '0.07746'
Set-TextValue 'E9' '  +4.34%  '
Set-TextValue 'D10' '1.123'
Set-TextValue 'E10' '  +2.86%  '
Set-TextValue 'D11' '42.01'
Set-TextValue 'E11' '  +0.27%  '
Set-TextValue 'E12' '  +3.59%  '
Set-TextValue 'D13' '6.349'
Set-TextValue 'E13' '  +3.75%  '
Set-TextValue 'B14' 'BinanceUSD'
Set-TextValue 'C14' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D14' '1.002'
Set-TextValue 'E14' '  +0.08%  '
Set-TextValue 'B15' 'Chainlink'
Set-TextValue 'C15' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D15' '7.625'
Set-TextValue 'E15' '  +5.34%  '
Set-TextValue 'D16' '1.828.75'
Set-TextValue 'E16' '  +2.46%  '
Set-TextValue 'E17' '  +3.18%  '
Set-TextValue 'D18' '89.97'
Set-TextValue 'E18' '  +1.10%  '
Set-TextValue 'D19' '0.06589'
Set-TextValue 'E19' '  +1.41%  '
Set-TextValue 'D20' '17.76'
Set-TextValue 'E20' '  +2.91%  '
Set-TextValue 'E21' '  +0.04%  '
Set-TextValue 'D22' '6.070'
Set-TextValue 'E22' '  +3.09%  '
Set-TextValue 'D23' '28.569.62'
Set-TextValue 'E23' '  +1.81%  '
Set-TextValue 'D24' '11.22'
Set-TextValue 'E24' '  +0.62%  '
Set-TextValue 'D25' '2.262'
Set-TextValue 'E25' '  +8.16%  '
Set-TextValue 'D26' '158.21'
Set-TextValue 'E26' '  +1.79%  '
Set-TextValue 'E27' '  +2.72%  '
Set-TextValue 'D28' '2.460'
Set-TextValue 'E28' '  +6.89%  '
Set-TextValue 'D29' '2.039.37'
Set-TextValue 'E29' '  +2.27%  '
Set-TextValue 'E30' '  +2.70%  '
Set-TextValue 'D31' '1.139'
Set-TextValue 'E31' '  +2.09%  '
Set-TextValue 'D32' '0.1123'
Set-TextValue 'E32' '  +6.20%  '
Set-TextValue 'D33' '5.702'
Set-TextValue 'E33' '  +2.63%  '
Set-TextValue 'D34' '0.07383'
Set-TextValue 'E34' '  +13.89%  '
Set-TextValue 'E35' '  -0.44%  '
Set-TextValue 'D36' '0.2261'
Set-TextValue 'E36' '  +0.64%  '
Set-TextValue 'D37' '0.02356'
Set-TextValue 'E37' '  +2.91%  '
Set-TextValue 'D38' '8.965'
Set-TextValue 'E38' '  +6.05%  '
Set-TextValue 'D39' '5.220'
Set-TextValue 'E39' '  +4.08%  '
Set-TextValue 'E40' '  +2.65%  '
Set-TextValue 'D41' '0.6304'
Set-TextValue 'E41' '  +1.88%  '
Set-TextValue 'D42' '1.192'
Set-TextValue 'E42' '  +1.07%  '
Set-TextValue 'E43' '  +0.09%  '
Set-TextValue 'D44' '1.400'
Set-TextValue 'E44' '  -3.43%  '
Set-TextValue 'E45' '  +1.59%  '
Set-TextValue 'D46' '0.5912'
Set-TextValue 'E46' '  +2.33%  '
Set-TextValue 'D47' '3.711'
Set-TextValue 'E47' '  +1.15%  '
Set-TextValue 'D48' '125.23'
Set-TextValue 'E48' '  +0.31%  '
Set-TextValue 'D49' '2.004'
Set-TextValue 'E49' '  +4.10%  '
Set-TextValue 'D50' '1.197'
Set-TextValue 'E50' '  +0.78%  '
Set-TextValue 'D51' '0.06921'
